# Add class quest manager.
# Introduces a new "Start_Quest_Dragon_king" dialogue/line/step used to
# kick off the Dragon King quest, and normalizes a few Dialogue "Type"
# values (Default -> Normal / Nomarl) on the Dialogues sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Dialogues sheet
# ---------------------------------------------------------------
$dialogues = $wb.Worksheets.Item("Dialogues")

$dialogues.Range("B3").Value = "Normal"
$dialogues.Range("B4").Value = "Nomarl"
$dialogues.Range("B5").Value = "Nomarl"

$dialogues.Range("A6").Value = "Start_Quest_Dragon_king"
$dialogues.Range("A6").HorizontalAlignment = -4108
$dialogues.Range("B6").Value = "Start"
$dialogues.Range("B6").HorizontalAlignment = -4108
$dialogues.Range("C6").Value = "Dragon_king_Eatern_Sea"

$dialogues.Columns.Item(1).ColumnWidth = 25.666666666666668

# ---------------------------------------------------------------
# Lines sheet
# ---------------------------------------------------------------
$lines = $wb.Worksheets.Item("Lines")

$lines.Range("B2").Value = "Start_Quest_Dragon_king"

$lines.Range("A17").Value = "Default-Dragon_king"
$lines.Range("B17").Value = "Default_Dragon_king"
$lines.Range("B17").HorizontalAlignment = -4108
$lines.Range("C17").Value = "Dragon_king_Eatern_Sea"
$lines.Range("D17").Value = "L_Deafult_Dragon_king_Eatern_Sea"
$lines.Range("E17").Value = $false

# ---------------------------------------------------------------
# Quests sheet - selection moves, no data changes
# ---------------------------------------------------------------
$quests = $wb.Worksheets.Item("Quests")

# ---------------------------------------------------------------
# Steps sheet
# ---------------------------------------------------------------
$steps = $wb.Worksheets.Item("Steps")

$steps.Range("D2").Value = "Start_Quest_Dragon_king"
$steps.Columns.Item(4).ColumnWidth = 20.666666666666668

# ---------------------------------------------------------------
# Rewards sheet - loses tabSelected (becomes inactive), no data changes
# ---------------------------------------------------------------
$rewards = $wb.Worksheets.Item("Rewards")

# ---------------------------------------------------------------
# Selection / active-sheet bookkeeping (order matters: the last
# Select() wins and becomes the active / tabSelected sheet).
# ---------------------------------------------------------------
$dialogues.Range("A6").Select()
$quests.Range("D30").Select()
$steps.Range("E3").Select()
